# Update the "取得日時" (retrieved datetime) column for all data rows on the
# "ランサーズ" sheet from 2025-10-06 01:17:09 to 2025-10-06 01:43:12, matching
# commit "Append: 2025-10-06 01:43 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-10-06 01:17:09"
$newTimestamp = "2025-10-06 01:43:12"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
